# Generate Report for Archive
# Status for c11fde86-...md (row 8) and c3b51552-...md (row 9) moves from
# "Ready for handoff" to "In Translation" across the Overview, zh-cn, and
# de-de sheets. Row 10 (e654846b-...md) stays "Ready for handoff".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = "In Translation"
$overview.Range("C8").Value = "In Translation"
$overview.Range("B9").Value = "In Translation"
$overview.Range("C9").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C8").Value = "In Translation"
$zhcn.Range("C9").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C8").Value = "In Translation"
$dede.Range("C9").Value = "In Translation"
